$wb = $excel.ActiveWorkbook

# Mapping of sheet name -> { row -> new F value }
$updates = @{
    "展览" = @{
        2  = 14825
        3  = 18368
        5  = 107
        14 = 93
        15 = 193
        16 = 52
        17 = 1397
        18 = 152
        22 = 7606
        24 = 17
        26 = 1209
        28 = 5934
        29 = 93
        30 = 58
        31 = 155
        34 = 5261
    }
    "全部类型" = @{
        2  = 14825
        3  = 18368
        5  = 107
        14 = 93
        15 = 193
        16 = 52
        17 = 1397
        18 = 152
        23 = 7606
        25 = 17
        27 = 1209
        30 = 5934
        31 = 93
        32 = 58
        33 = 155
        36 = 5261
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($r in $rows.Keys) {
        $ws.Cells.Item($r, 6).Value = $rows[$r]
    }
}
